$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1795.3721  # H15 was 1804.6976
$ws.Cells.Item(15, 9).Value = 1795.3721  # I15 was 1804.6976
$ws.Cells.Item(15, 11).Value = 5386.1163  # K15 was 5414.0928
$ws.Cells.Item(15, 13).Value = -5217.1163  # M15 was -5245.0928
$ws.Cells.Item(32, 8).Value = 1444.1111  # H32 was 1828.2858
$ws.Cells.Item(32, 9).Value = 1000.75  # I32 was 1239.8
$ws.Cells.Item(32, 10).Value = 1798.8  # J32 was 3299.5
$ws.Cells.Item(32, 11).Value = 1000.75  # K32 was 1239.8
$ws.Cells.Item(32, 12).Value = 1798.8  # L32 was 3299.5
$ws.Cells.Item(32, 13).Value = -674.75  # M32 was -913.8
$ws.Cells.Item(32, 14).Value = -2450.8  # N32 was -3951.5
$ws.Cells.Item(70, 8).Value = 33073046  # H70 was 31569794
$ws.Cells.Item(70, 9).Value = 50001308  # I70 was 41667944
$ws.Cells.Item(70, 10).Value = 27782964  # J70 was 27782988
$ws.Cells.Item(70, 11).Value = 150003924  # K70 was 125003832
$ws.Cells.Item(70, 12).Value = 83348892  # L70 was 83348964
$ws.Cells.Item(70, 13).Value = -150003654  # M70 was -125003562
$ws.Cells.Item(70, 14).Value = -83349432  # N70 was -83349504
$ws.Cells.Item(73, 8).Value = 33073046  # H73 was 31569794
$ws.Cells.Item(73, 9).Value = 50001308  # I73 was 41667944
$ws.Cells.Item(73, 10).Value = 27782964  # J73 was 27782988
$ws.Cells.Item(73, 11).Value = 150003924  # K73 was 125003832
$ws.Cells.Item(73, 12).Value = 83348892  # L73 was 83348964
$ws.Cells.Item(73, 13).Value = -150002988  # M73 was -125002896
$ws.Cells.Item(73, 14).Value = -83350764  # N73 was -83350836
$ws.Cells.Item(98, 8).Value = 3685.4814  # H98 was 3916.2173
$ws.Cells.Item(98, 9).Value = 3793.1155  # I98 was 4053.9092
$ws.Cells.Item(98, 11).Value = 3793.1155  # K98 was 4053.9092
$ws.Cells.Item(98, 13).Value = -2295.1155  # M98 was -2555.9092
$ws.Cells.Item(122, 8).Value = 3685.4814  # H122 was 3916.2173
$ws.Cells.Item(122, 9).Value = 3793.1155  # I122 was 4053.9092
$ws.Cells.Item(122, 11).Value = 11379.3465  # K122 was 12161.7276
$ws.Cells.Item(122, 13).Value = -8929.3465  # M122 was -9711.7276
$ws.Cells.Item(126, 8).Value = 87499.5  # H126 was 58569.332
$ws.Cells.Item(126, 9).Value = 0  # I126 was 709
$ws.Cells.Item(126, 11).Value = 0  # K126 was 709
$ws.Cells.Item(126, 13).ClearContents()  # M126 was 4231
$ws.Cells.Item(127, 8).Value = 1000  # H127 was 763.8570999999999
$ws.Cells.Item(127, 9).Value = 1000  # I127 was 763.8570999999999
$ws.Cells.Item(127, 11).Value = 3000  # K127 was 2291.5713
$ws.Cells.Item(127, 13).Value = 1960  # M127 was 2668.4287
$ws.Cells.Item(132, 8).Value = 1305.4445  # H132 was 1320.174
$ws.Cells.Item(132, 10).Value = 2659.8333  # J132 was 2563.1428
$ws.Cells.Item(132, 12).Value = 7979.499899999999  # L132 was 7689.428400000001
$ws.Cells.Item(132, 14).Value = -13039.4999  # N132 was -12749.4284
$ws.Cells.Item(135, 8).Value = 6667694.5  # H135 was 7693439.5
$ws.Cells.Item(135, 9).Value = 7693033  # I135 was 9091701
$ws.Cells.Item(135, 11).Value = 69237297  # K135 was 81825309
$ws.Cells.Item(135, 13).Value = -69234762  # M135 was -81822774
$ws.Cells.Item(137, 8).Value = 5659.0195  # H137 was 5659.039
$ws.Cells.Item(137, 9).Value = 4224.7188  # I137 was 4126.8486
$ws.Cells.Item(137, 10).Value = 8074.684  # J137 was 8468.056
$ws.Cells.Item(137, 11).Value = 12674.1564  # K137 was 12380.5458
$ws.Cells.Item(137, 12).Value = 24224.052  # L137 was 25404.168
$ws.Cells.Item(137, 13).Value = -10124.1564  # M137 was -9830.5458
$ws.Cells.Item(137, 14).Value = -29324.052  # N137 was -30504.168
$ws.Cells.Item(138, 8).Value = 4098.826  # H138 was 4028.4167
$ws.Cells.Item(138, 10).Value = 4751.4707  # J138 was 4621.3335
$ws.Cells.Item(138, 12).Value = 14254.4121  # L138 was 13864.0005
$ws.Cells.Item(138, 14).Value = -24534.4121  # N138 was -24144.0005
$ws.Cells.Item(141, 8).Value = 3126.4285  # H141 was 2821.875
$ws.Cells.Item(141, 9).Value = 3126.4285  # I141 was 2821.875
$ws.Cells.Item(141, 11).Value = 9379.2855  # K141 was 8465.625
$ws.Cells.Item(141, 13).Value = -4199.2855  # M141 was -3285.625
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 29209.105  # H74 was 28485.795
$ws.Cells.Item(74, 10).Value = 4699.067  # J74 was 4467.875
$ws.Cells.Item(74, 12).Value = 4699.067  # L74 was 4467.875
$ws.Cells.Item(74, 14).Value = -6447.067  # N74 was -6215.875
$ws.Cells.Item(77, 8).Value = 29209.105  # H77 was 28485.795
$ws.Cells.Item(77, 10).Value = 4699.067  # J77 was 4467.875
$ws.Cells.Item(77, 12).Value = 23495.335  # L77 was 22339.375
$ws.Cells.Item(77, 14).Value = -32231.335  # N77 was -31075.375
$ws.Cells.Item(97, 8).Value = 2689151.2  # H97 was 2977253.2
$ws.Cells.Item(97, 9).Value = 859.1429000000001  # I97 was 938.28
$ws.Cells.Item(97, 11).Value = 859.1429000000001  # K97 was 938.28
$ws.Cells.Item(97, 13).Value = -363.1429000000001  # M97 was -442.28
$ws.Cells.Item(122, 8).Value = 3524.4092  # H122 was 3596.9524
$ws.Cells.Item(122, 9).Value = 2752.7334  # I122 was 2806.4285
$ws.Cells.Item(122, 11).Value = 8258.200199999999  # K122 was 8419.2855
$ws.Cells.Item(122, 13).Value = -5808.200199999999  # M122 was -5969.2855
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(26, 8).Value = 6832.6665  # H26 was 5289.4
$ws.Cells.Item(26, 9).Value = 6832.6665  # I26 was 5289.4
$ws.Cells.Item(26, 11).Value = 6832.6665  # K26 was 5289.4
$ws.Cells.Item(26, 13).Value = -6540.6665  # M26 was -4997.4
$ws.Cells.Item(28, 8).Value = 53842  # H28 was 41921
$ws.Cells.Item(28, 10).Value = 53842  # J28 was 41921
$ws.Cells.Item(28, 12).Value = 53842  # L28 was 41921
$ws.Cells.Item(28, 14).Value = -54430  # N28 was -42509
$ws.Cells.Item(68, 8).Value = 43000  # H68 was 71963.336
$ws.Cells.Item(68, 10).Value = 43000  # J68 was 71963.336
$ws.Cells.Item(68, 12).Value = 43000  # L68 was 71963.336
$ws.Cells.Item(68, 14).Value = -44622  # N68 was -73585.336
$ws.Cells.Item(71, 8).Value = 43000  # H71 was 71963.336
$ws.Cells.Item(71, 10).Value = 43000  # J71 was 71963.336
$ws.Cells.Item(71, 12).Value = 129000  # L71 was 215890.008
$ws.Cells.Item(71, 14).Value = -137112  # N71 was -224002.008
$ws.Cells.Item(86, 8).Value = 89842.75  # H86 was 106900.1
$ws.Cells.Item(86, 9).Value = 148997.86  # I86 was 172999.17
$ws.Cells.Item(86, 10).Value = 7025.6  # J86 was 7751.5
$ws.Cells.Item(86, 11).Value = 148997.86  # K86 was 172999.17
$ws.Cells.Item(86, 12).Value = 7025.6  # L86 was 7751.5
$ws.Cells.Item(86, 13).Value = -147874.86  # M86 was -171876.17
$ws.Cells.Item(86, 14).Value = -9271.6  # N86 was -9997.5
$ws.Cells.Item(89, 8).Value = 89842.75  # H89 was 106900.1
$ws.Cells.Item(89, 9).Value = 148997.86  # I89 was 172999.17
$ws.Cells.Item(89, 10).Value = 7025.6  # J89 was 7751.5
$ws.Cells.Item(89, 11).Value = 744989.2999999999  # K89 was 864995.8500000001
$ws.Cells.Item(89, 12).Value = 35128  # L89 was 38757.5
$ws.Cells.Item(89, 13).Value = -739373.2999999999  # M89 was -859379.8500000001
$ws.Cells.Item(89, 14).Value = -46360  # N89 was -49989.5
$ws.Cells.Item(94, 8).Value = 3616.3635  # H94 was 3779.0476
$ws.Cells.Item(94, 9).Value = 1727.5834  # I94 was 1866.4546
$ws.Cells.Item(94, 11).Value = 1727.5834  # K94 was 1866.4546
$ws.Cells.Item(94, 13).Value = -1276.5834  # M94 was -1415.4546
$ws.Cells.Item(134, 8).Value = 4171208.2  # H134 was 4102838.5
$ws.Cells.Item(134, 9).Value = 6099533.5  # I134 was 5954322
$ws.Cells.Item(134, 11).Value = 18298600.5  # K134 was 17862966
$ws.Cells.Item(134, 13).Value = -18296065.5  # M134 was -17860431
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 8610.647000000001  # H31 was 7923.6313
$ws.Cells.Item(31, 9).Value = 2546.3076  # I31 was 2402.4375
$ws.Cells.Item(31, 10).Value = 12364.762  # J31 was 11939.046
$ws.Cells.Item(31, 11).Value = 2546.3076  # K31 was 2402.4375
$ws.Cells.Item(31, 12).Value = 12364.762  # L31 was 11939.046
$ws.Cells.Item(31, 13).Value = -2251.3076  # M31 was -2107.4375
$ws.Cells.Item(31, 14).Value = -12954.762  # N31 was -12529.046
$ws.Cells.Item(34, 8).Value = 8610.647000000001  # H34 was 7923.6313
$ws.Cells.Item(34, 9).Value = 2546.3076  # I34 was 2402.4375
$ws.Cells.Item(34, 10).Value = 12364.762  # J34 was 11939.046
$ws.Cells.Item(34, 11).Value = 2546.3076  # K34 was 2402.4375
$ws.Cells.Item(34, 12).Value = 12364.762  # L34 was 11939.046
$ws.Cells.Item(34, 13).Value = -2344.3076  # M34 was -2200.4375
$ws.Cells.Item(34, 14).Value = -12768.762  # N34 was -12343.046
$ws.Cells.Item(37, 8).Value = 10000  # H37 was 5500
$ws.Cells.Item(37, 9).Value = 0  # I37 was 1000
$ws.Cells.Item(37, 11).Value = 0  # K37 was 1000
$ws.Cells.Item(37, 13).ClearContents()  # M37 was -893
$ws.Cells.Item(62, 8).Value = 9222.727999999999  # H62 was 9280
$ws.Cells.Item(62, 9).Value = 8941.5  # I62 was 8999.799999999999
$ws.Cells.Item(62, 11).Value = 8941.5  # K62 was 8999.799999999999
$ws.Cells.Item(62, 13).Value = -8317.5  # M62 was -8375.799999999999
$ws.Cells.Item(65, 8).Value = 9222.727999999999  # H65 was 9280
$ws.Cells.Item(65, 9).Value = 8941.5  # I65 was 8999.799999999999
$ws.Cells.Item(65, 11).Value = 44707.5  # K65 was 44999
$ws.Cells.Item(65, 13).Value = -41587.5  # M65 was -41879
$ws.Cells.Item(105, 8).Value = 5103250  # H105 was 5495750
$ws.Cells.Item(105, 9).Value = 5953167.5  # I105 was 6494296
$ws.Cells.Item(105, 11).Value = 5953167.5  # K105 was 6494296
$ws.Cells.Item(105, 13).Value = -5951420.5  # M105 was -6492549
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 3238.2307  # H113 was 3386.25
$ws.Cells.Item(113, 9).Value = 1198  # I113 was 1168
$ws.Cells.Item(113, 10).Value = 3408.25  # J113 was 3829.9
$ws.Cells.Item(113, 11).Value = 3594  # K113 was 3504
$ws.Cells.Item(113, 12).Value = 10224.75  # L113 was 11489.7
$ws.Cells.Item(113, 13).Value = -1424  # M113 was -1334
$ws.Cells.Item(113, 14).Value = -14564.75  # N113 was -15829.7
$ws.Cells.Item(122, 8).Value = 1348688.9  # H122 was 1490613.5
$ws.Cells.Item(122, 9).Value = 1769214.8  # I122 was 2177380
$ws.Cells.Item(122, 10).Value = 3006  # J122 was 2619.6667
$ws.Cells.Item(122, 11).Value = 15922933.2  # K122 was 19596420
$ws.Cells.Item(122, 12).Value = 27054  # L122 was 23577.0003
$ws.Cells.Item(122, 13).Value = -15920483.2  # M122 was -19593970
$ws.Cells.Item(122, 14).Value = -31954  # N122 was -28477.0003
$ws.Cells.Item(128, 8).Value = 120563.43  # H128 was 144243
$ws.Cells.Item(128, 9).Value = 120563.43  # I128 was 144243
$ws.Cells.Item(128, 11).Value = 361690.29  # K128 was 432729
$ws.Cells.Item(128, 13).Value = -356710.29  # M128 was -427749
$ws.Cells.Item(131, 8).Value = 37536.57  # H131 was 36307.277
$ws.Cells.Item(131, 10).Value = 49432.668  # J131 was 47271.5
$ws.Cells.Item(131, 12).Value = 148298.004  # L131 was 141814.5
$ws.Cells.Item(131, 14).Value = -158378.004  # N131 was -151894.5
$ws.Cells.Item(132, 8).Value = 4451.311  # H132 was 4626.905
$ws.Cells.Item(132, 9).Value = 3046.4285  # I132 was 3099.5
$ws.Cells.Item(132, 10).Value = 5680.5835  # J132 was 6015.4546
$ws.Cells.Item(132, 11).Value = 27417.8565  # K132 was 27895.5
$ws.Cells.Item(132, 12).Value = 51125.2515  # L132 was 54139.0914
$ws.Cells.Item(132, 13).Value = -24887.8565  # M132 was -25365.5
$ws.Cells.Item(132, 14).Value = -56185.2515  # N132 was -59199.0914
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 65121.53  # H80 was 61620.332
$ws.Cells.Item(80, 9).Value = 7098  # I80 was 6681.5
$ws.Cells.Item(80, 11).Value = 7098  # K80 was 6681.5
$ws.Cells.Item(80, 13).Value = -6100  # M80 was -5683.5
$ws.Cells.Item(83, 8).Value = 65121.53  # H83 was 61620.332
$ws.Cells.Item(83, 9).Value = 7098  # I83 was 6681.5
$ws.Cells.Item(83, 11).Value = 35490  # K83 was 33407.5
$ws.Cells.Item(83, 13).Value = -30498  # M83 was -28415.5
$ws.Cells.Item(126, 8).Value = 7006.6274  # H126 was 7016.451
$ws.Cells.Item(126, 9).Value = 4092.75  # I126 was 4113.5835
$ws.Cells.Item(126, 10).Value = 9596.741  # J126 was 9596.777
$ws.Cells.Item(126, 11).Value = 12278.25  # K126 was 12340.7505
$ws.Cells.Item(126, 12).Value = 28790.223  # L126 was 28790.331
$ws.Cells.Item(126, 13).Value = -9808.25  # M126 was -9870.750499999998
$ws.Cells.Item(126, 14).Value = -33730.223  # N126 was -33730.331
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4588.811  # H40 was 4685.7026
$ws.Cells.Item(40, 9).Value = 4125.9  # I40 was 4173.129
$ws.Cells.Item(40, 10).Value = 6572.7144  # J40 was 7334
$ws.Cells.Item(40, 11).Value = 4125.9  # K40 was 4173.129
$ws.Cells.Item(40, 12).Value = 6572.7144  # L40 was 7334
$ws.Cells.Item(40, 13).Value = -3989.9  # M40 was -4037.129
$ws.Cells.Item(40, 14).Value = -6844.7144  # N40 was -7606
$ws.Cells.Item(68, 8).Value = 5420.3125  # H68 was 4956.4
$ws.Cells.Item(68, 9).Value = 3873.375  # I68 was 3787.7778
$ws.Cells.Item(68, 10).Value = 6967.25  # J68 was 5912.5454
$ws.Cells.Item(68, 11).Value = 3873.375  # K68 was 3787.7778
$ws.Cells.Item(68, 12).Value = 6967.25  # L68 was 5912.5454
$ws.Cells.Item(68, 13).Value = -3124.375  # M68 was -3038.7778
$ws.Cells.Item(68, 14).Value = -8465.25  # N68 was -7410.5454
$ws.Cells.Item(71, 8).Value = 5420.3125  # H71 was 4956.4
$ws.Cells.Item(71, 9).Value = 3873.375  # I71 was 3787.7778
$ws.Cells.Item(71, 10).Value = 6967.25  # J71 was 5912.5454
$ws.Cells.Item(71, 11).Value = 19366.875  # K71 was 18938.889
$ws.Cells.Item(71, 12).Value = 34836.25  # L71 was 29562.727
$ws.Cells.Item(71, 13).Value = -15622.875  # M71 was -15194.889
$ws.Cells.Item(71, 14).Value = -42324.25  # N71 was -37050.727
$ws.Cells.Item(82, 8).Value = 2769.111  # H82 was 2996.625
$ws.Cells.Item(82, 9).Value = 1349.75  # I82 was 1483.3334
$ws.Cells.Item(82, 11).Value = 1349.75  # K82 was 1483.3334
$ws.Cells.Item(82, 13).Value = -988.75  # M82 was -1122.3334
$ws.Cells.Item(85, 8).Value = 2769.111  # H85 was 2996.625
$ws.Cells.Item(85, 9).Value = 1349.75  # I85 was 1483.3334
$ws.Cells.Item(85, 11).Value = 1349.75  # K85 was 1483.3334
$ws.Cells.Item(85, 13).Value = -101.75  # M85 was -235.3334
$ws.Cells.Item(104, 8).Value = 27428.334  # H104 was 35567
$ws.Cells.Item(104, 10).Value = 27428.334  # J104 was 35567
$ws.Cells.Item(104, 12).Value = 27428.334  # L104 was 35567
$ws.Cells.Item(104, 14).Value = -34416.334  # N104 was -42555
$ws.Cells.Item(122, 8).Value = 3769.5625  # H122 was 3682.9395
$ws.Cells.Item(122, 9).Value = 3012.84  # I122 was 3013.32
$ws.Cells.Item(122, 10).Value = 6472.143  # J122 was 5775.5
$ws.Cells.Item(122, 11).Value = 9038.52  # K122 was 9039.960000000001
$ws.Cells.Item(122, 12).Value = 19416.429  # L122 was 17326.5
$ws.Cells.Item(122, 13).Value = -6588.52  # M122 was -6589.960000000001
$ws.Cells.Item(122, 14).Value = -24316.429  # N122 was -22226.5
$ws.Cells.Item(132, 8).Value = 15165740  # H132 was 14299299
$ws.Cells.Item(132, 9).Value = 35719228  # I132 was 33338084
$ws.Cells.Item(132, 10).Value = 21064.79  # J132 was 20211.55
$ws.Cells.Item(132, 11).Value = 107157684  # K132 was 100014252
$ws.Cells.Item(132, 12).Value = 63194.37  # L132 was 60634.64999999999
$ws.Cells.Item(132, 13).Value = -107155154  # M132 was -100011722
$ws.Cells.Item(132, 14).Value = -68254.37  # N132 was -65694.64999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 6176.0884  # H132 was 6323.9395
$ws.Cells.Item(132, 9).Value = 6124.609  # I132 was 6344.0454
$ws.Cells.Item(132, 11).Value = 18373.827  # K132 was 19032.1362
$ws.Cells.Item(132, 13).Value = -15843.827  # M132 was -16502.1362

Write-Output "Done applying 243 cell updates across 8 sheets"